$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (J6/K6): "Reason of the Change" / "Changes Done" text updated for the
# hana "or"/"and" + CI-vs-DB control-flow fix.
$ws.Range("J6").Value = "1.If condition to check whether the DB is hana had an ""or"" condition.SO, when any of the condition was satisfied it got into the loop.`n2. The script has to run separately for CI and DB, but the both CI and DB are getting executed in the same run."
$ws.Range("K6").Value = "1.The control needs to get into the loop only when all the conditions are satisfied, so the ""or"" condition is changed to ""and"" condition. This ensures that the control does not get into the loop even if one of the condition fails.`n2. To make the CI and DB part to execute in separate run, the ""if"" condition of DB is changed to ""elif"", so that in one run only one of the condition gets satisfied."

# Row 8 (K8): "Changes Done" text updated to explain the sidadm/hdbsid fix in
# more detail; the wrapped cell grows taller, so the row auto-sizes to 43.5.
$ws.Range("K8").Value = "1.The DB user that was present was ""sidadm"" not  ""hdbsid(like in oracle)"". So, the db user is changed to ""sidadm""."
$ws.Rows.Item(8).RowHeight = 43.5

# Selection moved to J8.
$ws.Range("J8").Select() | Out-Null
